# cs-en-us-pbqs.xlsx weekly refresh: new crime data collected.
# Updates the report header (volume/number + week-covering dates) and the
# CompStat weekly/28-day/YTD crime-count table (rows 14-30) to the new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume number and week-covering date range -----------
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Crime-count table (columns: C/D/E = Week-to-Date 23/22/%chg,
#     F/G/H = 28-Day 23/22/%chg, I/J/K = YTD 23/22/%chg,
#     L/M/N = 2yr/13yr/30yr %chg) ---------------------------------------------
# Row 14 - Murder
$ws.Range("G14").Value = 2

# Row 15 - Rape
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -6.666666666666
$ws.Range("L15").Value = 27.272727272727
$ws.Range("M15").Value = 27.272727272727
$ws.Range("N15").Value = -56.25

# Row 16 - Robbery
$ws.Range("C16").Value = 23
$ws.Range("D16").Value = 23
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 76
$ws.Range("G16").Value = 112
$ws.Range("H16").Value = -32.142857142857
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 139
$ws.Range("K16").Value = -28.776978417266
$ws.Range("L16").Value = 7.608695652173
$ws.Range("M16").Value = -61.023622047244
$ws.Range("N16").Value = -88.228299643281

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 45
$ws.Range("D17").Value = 42
$ws.Range("E17").Value = 7.142857142857
$ws.Range("F17").Value = 185
$ws.Range("G17").Value = 162
$ws.Range("H17").Value = 14.197530864197
$ws.Range("I17").Value = 249
$ws.Range("J17").Value = 219
$ws.Range("K17").Value = 13.698630136986
$ws.Range("L17").Value = 25.125628140703
$ws.Range("M17").Value = 81.751824817518
$ws.Range("N17").Value = -25.225225225225

# Row 18 - Burglary
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = 21.428571428571
$ws.Range("F18").Value = 71
$ws.Range("G18").Value = 74
$ws.Range("H18").Value = -4.054054054054
$ws.Range("I18").Value = 103
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 19.767441860465
$ws.Range("M18").Value = -48.756218905472
$ws.Range("N18").Value = -89.192025183630

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 66
$ws.Range("D19").Value = 68
$ws.Range("E19").Value = -2.941176470588
$ws.Range("F19").Value = 253
$ws.Range("G19").Value = 292
$ws.Range("H19").Value = -13.356164383561
$ws.Range("I19").Value = 329
$ws.Range("J19").Value = 372
$ws.Range("K19").Value = -11.559139784946
$ws.Range("L19").Value = 55.924170616113
$ws.Range("M19").Value = 20.512820512820
$ws.Range("N19").Value = -39.633027522935

# Row 20 - G.L.A.
$ws.Range("C20").Value = 25
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = 4.166666666666
$ws.Range("F20").Value = 122
$ws.Range("G20").Value = 139
$ws.Range("H20").Value = -12.230215827338
$ws.Range("I20").Value = 163
$ws.Range("J20").Value = 165
$ws.Range("K20").Value = -1.212121212121
$ws.Range("L20").Value = 64.646464646464
$ws.Range("M20").Value = -11.413043478260
$ws.Range("N20").Value = -91.85

# Row 21 - TOTAL
$ws.Range("C21").Value = 178
$ws.Range("D21").Value = 173
$ws.Range("E21").Value = 2.890173410404
$ws.Range("F21").Value = 719
$ws.Range("G21").Value = 791
$ws.Range("H21").Value = -9.102402022756
$ws.Range("I21").Value = 957
$ws.Range("J21").Value = 1016
$ws.Range("K21").Value = -5.807086614173
$ws.Range("L21").Value = 36.714285714285
$ws.Range("M21").Value = -10.056390977443
$ws.Range("N21").Value = -79.711681153275

# Row 22 - Transit
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 400
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 20
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = 15.384615384615
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = 66.666666666666

# Row 23 - Housing
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 66.666666666666
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = 5.263157894736
$ws.Range("L23").Value = 11.111111111111
$ws.Range("M23").Value = -4.761904761904

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 160
$ws.Range("D24").Value = 129
$ws.Range("E24").Value = 24.031007751938
$ws.Range("F24").Value = 721
$ws.Range("G24").Value = 634
$ws.Range("H24").Value = 13.722397476340
$ws.Range("I24").Value = 867
$ws.Range("J24").Value = 784
$ws.Range("K24").Value = 10.586734693877
$ws.Range("L24").Value = 25.470332850940
$ws.Range("M24").Value = 50

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 62
$ws.Range("D25").Value = 72
$ws.Range("E25").Value = -13.888888888888
$ws.Range("F25").Value = 301
$ws.Range("G25").Value = 265
$ws.Range("H25").Value = 13.584905660377
$ws.Range("I25").Value = 394
$ws.Range("J25").Value = 329
$ws.Range("K25").Value = 19.756838905775
$ws.Range("L25").Value = 55.731225296442
$ws.Range("M25").Value = -8.584686774941

# Row 26 - UCR Rape*
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 22
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 15.789473684210

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -25
$ws.Range("F27").Value = 36
$ws.Range("G27").Value = 26
$ws.Range("H27").Value = 38.461538461538
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 31.25
$ws.Range("L27").Value = 27.272727272727

# Row 28 - Shooting Vic.
$ws.Range("C28").Value = 2
$ws.Range("C28").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 140
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = 133.333333333333
$ws.Range("L28").Value = 27.272727272727
$ws.Range("M28").Value = 55.555555555555
$ws.Range("N28").Value = -76.271186440678

# Row 29 - Shooting Inc.
$ws.Range("C29").Value = 2
$ws.Range("C29").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 60
$ws.Range("I29").Value = 9
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = 12.5
$ws.Range("M29").Value = 12.5
$ws.Range("N29").Value = -83.333333333333

# Row 30 - Hate Crimes
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = 300
